$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "Last updated" date from 2024-03-15 to 2024-03-28 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- "RAF-capacity" sheet: raise the capacity-credit multiplier for the two ---
# --- hydrogen technologies (rows 24 & 25, column B) from 0.3 to 1           ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Slightly narrow column A on the capacity sheet (user resized it).
$wsCapacity.Columns.Item(1).ColumnWidth = 28.17

# Scroll/zoom state + selected cell on the capacity sheet.
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 80

# The capacity sheet becomes the active / selected tab (RAF-generation loses it).
$wsCapacity.Select()
